# Swap the "B" and "C" quarter rows (2nd and 3rd rows) within each year
# block, then delete the now-duplicated/derived columns F (产销率) and G
# (销售量), which are dropped from the sheet entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $cols = @("A", "B", "C", "D", "E")
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

# Year blocks are laid out as rows (A,B,C,D) -> e.g. 2016: rows 2-5,
# 2017: rows 6-9, 2018: rows 10-13, 2019: rows 14-17.
# The "B" row and "C" row of each block (2nd and 3rd rows) are swapped.
Swap-Rows 3 4
Swap-Rows 7 8
Swap-Rows 11 12
Swap-Rows 15 16

# Remove columns F and G (载货汽车产销率 / 载货汽车销售量) entirely.
$ws.Range("F1:G17").Delete()

Write-Output "done"
